{"js": "// Apply the documented edits to the Word \"Function Description\" spec document.\n//\n// 1. Fill in the function name placeholder.\n// 2. Append the parameter declaration after \"Parameter List:\".\n// 3. Fill in the first data row of the parameter table (Name / Type / Description).\n// 4. Fill in the \"Returns:\" placeholder.\n// 5. Replace the \"Description:\" placeholder with the real description text,\n//    and append the follow-up explanation paragraphs.\n\nconst body = context.document.body;\n\n// --- 1. \"Name of function\" -> \"validatePackageBox\" -----------------------\nconst nameResults = body.search(\"Name of function\", { matchCase: true });\nnameResults.load(\"items\");\nawait context.sync();\nnameResults.items[0].insertText(\"validatePackageBox\", Word.InsertLocation.replace);\n\n// --- 2. \"Parameter List:\" gets a new bold run \" double num\" --------------\nconst paramListResults = body.search(\"Parameter List:\", { matchCase: true });\nparamListResults.load(\"items\");\nawait context.sync();\nconst paramListRange = paramListResults.items[0];\nconst paramExtra = paramListRange.insertText(\" double num\", Word.InsertLocation.after);\nparamExtra.font.bold = true;\n\nawait context.sync();\n\n// --- 3. First data row of the parameter table -----------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst nameCell = table.getCell(1, 0);\nnameCell.body.paragraphs.load(\"items\");\nconst typeCell = table.getCell(1, 1);\ntypeCell.body.paragraphs.load(\"items\");\nconst descCell = table.getCell(1, 2);\ndescCell.body.paragraphs.load(\"items\");\nawait context.sync();\n\nnameCell.body.paragraphs.items[0].insertText(\"Num\", Word.InsertLocation.end);\ntypeCell.body.paragraphs.items[0].insertText(\"Double\", Word.InsertLocation.end);\ndescCell.body.paragraphs.items[0].insertText(\n  \"A double type value representing the user\\u2019s input for the package\\u2019s size.\",\n  Word.InsertLocation.end\n);\n\nawait context.sync();\n\n// --- 4. \"Returns:\" placeholder --------------------------------------------\nconst returnsResults = body.search(\n  \"Return type and description of what it means and special conditions that affect it.\",\n  { matchCase: true }\n);\nreturnsResults.load(\"items\");\nawait context.sync();\nreturnsResults.items[0].insertText(\"int, 1 = true and 0 = false\", Word.InsertLocation.replace);\n\n// --- 5. \"Description:\" placeholder + follow-up paragraphs ------------------\nconst descPlaceholder =\n  \"A description of what the function does, any special algorithms used and special condition that the user needs to be aware of that will affect the output. There needs to be sufficient detail in the description to allow the black box tests to be written before the code is complete. This description can also be given to the programmers and provide them with everything they need to know to write the code.\";\n\nconst descResults = body.search(descPlaceholder, { matchCase: true });\ndescResults.load(\"items\");\nawait context.sync();\n\nconst descRange = descResults.items[0];\nconst descNewText =\n  \"This function returns a valid value if the parameter is within the limitation as per project requirements. The macros defined in input.h is used: \" +\n  \"SZMIN, SZMID, SZMAX \" +\n  \"to perform a relational operation of \" +\n  \"||.\" +\n  \" \";\ndescRange.insertText(descNewText, Word.InsertLocation.replace);\n\nawait context.sync();\n\n// Paragraph holding \"Description:\" is the last paragraph in the body; add the\n// two new paragraphs after it (they inherit the grey \"placeholder\" character\n// formatting from the run they are split off from, matching the template).\nconst bodyParagraphs = body.paragraphs;\nbodyParagraphs.load(\"items\");\nawait context.sync();\n\nconst descParagraph = bodyParagraphs.items[bodyParagraphs.items.length - 1];\nconst ifParagraphText =\n  \"IF num \" +\n  \"is any of the values set in the macros above then \" +\n  \"1 is returned. Else, 0 is returned indicating an invalid entry regarding the package\\u2019s weight is entered.  \";\nconst ifParagraph = descParagraph.insertParagraph(ifParagraphText, Word.InsertLocation.after);\nifParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Apply the documented edits to the Word \"Function Description\" spec document.\n#\n# 1. Fill in the function name placeholder.\n# 2. Append the parameter declaration after \"Parameter List:\".\n# 3. Fill in the first data row of the parameter table (Name / Type / Description).\n# 4. Fill in the \"Returns:\" placeholder.\n# 5. Replace the \"Description:\" placeholder with the real description text,\n#    and append the follow-up explanation paragraphs.\n\n$d = $word.ActiveDocument\n\n# Unicode right single quotation mark used in \"user's\" / \"package's\".\n$rsquo = [char]0x2019\n\n# --- 1. \"Name of function\" -> \"validatePackageBox\" -------------------------\n$d.Content.Find.Execute(\"Name of function\", $false, $false, $false, $false, $false, $true, 1, $false, \"validatePackageBox\", 2) | Out-Null\n\n# --- 2. \"Parameter List:\" gets a new bold run \" double num\" ----------------\n$paramRange = $d.Content\n$paramRange.Find.Execute(\"Parameter List:\") | Out-Null\n$paramRange.Collapse(0)\n$paramRange.InsertAfter(\" double num\")\n$paramRange.Bold = 1\n\n# --- 3. First data row of the parameter table -------------------------------\n$table = $d.Tables.Item(1)\n$table.Cell(2, 1).Range.InsertBefore(\"Num\")\n$table.Cell(2, 2).Range.InsertBefore(\"Double\")\n$descText = \"A double type value representing the user\" + $rsquo + \"s input for the package\" + $rsquo + \"s size.\"\n$table.Cell(2, 3).Range.InsertBefore($descText)\n\n# --- 4. \"Returns:\" placeholder ----------------------------------------------\n$d.Content.Find.Execute(\"Return type and description of what it means and special conditions that affect it.\", $false, $false, $false, $false, $false, $true, 1, $false, \"int, 1 = true and 0 = false\", 2) | Out-Null\n\n# --- 5. \"Description:\" placeholder + follow-up paragraphs -------------------\n$oldDescription = \"A description of what the function does, any special algorithms used and special condition that the user needs to be aware of that will affect the output. There needs to be sufficient detail in the description to allow the black box tests to be written before the code is complete. This description can also be given to the programmers and provide them with everything they need to know to write the code.\"\n$newDescription = \"This function returns a valid value if the parameter is within the limitation as per project requirements. The macros defined in input.h is used: SZMIN, SZMID, SZMAX to perform a relational operation of ||. \"\n$d.Content.Find.Execute($oldDescription, $false, $false, $false, $false, $false, $true, 1, $false, $newDescription, 2) | Out-Null\n\n# The paragraph holding \"Description:\" is now the last paragraph in the body;\n# add the two new paragraphs after it. They inherit the grey placeholder\n# character formatting from the paragraph mark they are split off from.\n# NOTE: use $d.Content.Paragraphs.Last (not $d.Paragraphs.Last) \u2014 after the\n# table-cell edits above, $d.Paragraphs' own index cache goes stale.\n$descParagraph = $d.Content.Paragraphs.Last\n$descParagraph.Range.InsertParagraphAfter()\n\n$ifParagraph = $d.Content.Paragraphs.Last\n$ifText = \"IF num is any of the values set in the macros above then 1 is returned. Else, 0 is returned indicating an invalid entry regarding the package\" + $rsquo + \"s weight is entered.  \"\n$ifParagraph.Range.Text = $ifText\n\n$ifParagraph.Range.InsertParagraphAfter()\n\nWrite-Output \"done\"\n"}
